# Adding more student feedback: the "jobs" sheet gains new comparison rows
# (Engineering as a core focus / Business Problem Solving) while a few of
# the old rows (Engineering, Modern ML math/research, Problem solving using
# machine learning) are dropped, shrinking the table from 10 data rows
# (A2:D11) down to 7 data rows (A2:D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing three rows that no longer exist in the updated table.
# This also shrinks the sheet's used range/dimension down to row 8.
$ws.Range("A9:D11").Delete()

# The conditional formatting (green >0.1 / red <-0.1 on column D) should
# keep applying to exactly the new data extent.
$cfRange = $ws.Range("D2:D11")
$cfRange.FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D8"))

# Row 3 becomes "Engineering as a core focus"
$ws.Range("A3").Value = "Engineering as a core focus"
$ws.Range("B3").Value = 0.33
$ws.Range("C3").Value = 0.38
$ws.Range("D3").Value = -0.04

# Row 4 becomes "Teaches SQL/database languages"
$ws.Range("A4").Value = "Teaches SQL/database languages"
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0

# Row 5 becomes "Traditional Machine Learning math/statistics"
$ws.Range("A5").Value = "Traditional Machine Learning math/statistics"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0

# Row 6 becomes "Business Problem Solving"
$ws.Range("A6").Value = "Business Problem Solving"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.62
$ws.Range("D6").Value = 0.38

# Row 7 becomes "Teaches Big data technologies"
$ws.Range("A7").Value = "Teaches Big data technologies"
$ws.Range("B7").Value = 0.33
$ws.Range("C7").Value = 0.75
$ws.Range("D7").Value = -0.42

# Row 8 stays "Business Communication/ interaction" with updated numbers
$ws.Range("A8").Value = "Business Communication/ interaction"
$ws.Range("B8").Value = 0.83
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.33
